$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.920.68"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "'1.923.17"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'320.26"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'0.5060"
$ws.Range("E7").Value = "  -2.41%  "

$ws.Range("D8").Value = "'0.4035"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "'0.08343"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").Value = "'42.51"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").Value = "'1.104"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "'23.78"
$ws.Range("E12").Value = "  +2.30%  "

$ws.Range("D13").Value = "'1.918.65"
$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").Value = "'6.410"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "'7.237"
$ws.Range("E15").Value = "  -1.06%  "

$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "'92.23"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").Value = "'0.00001098"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").Value = "'0.06511"
$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").Value = "'18.28"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D22").Value = "'5.950"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "'29.936.36"

$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").Value = "'2.191"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  +2.37%  "

$ws.Range("D27").Value = "'2.141.57"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").Value = "'162.07"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").Value = "'2.325"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").Value = "'129.02"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").Value = "'1.131"
$ws.Range("E31").Value = "  +3.85%  "

$ws.Range("D32").Value = "'0.1037"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").Value = "'5.967"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").Value = "'3.789"
$ws.Range("E34").Value = "  +1.19%  "

$ws.Range("D35").Value = "'0.02452"
$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("D36").Value = "'5.402"
$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").Value = "'0.06422"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").Value = "'0.2157"
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").Value = "'0.6520"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "'8.740"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("D41").Value = "'1.195"
$ws.Range("E41").Value = "  -1.98%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.40"
$ws.Range("E42").Value = "  -3.57%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.219"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").Value = "'2.239"
$ws.Range("E44").Value = "  +8.93%  "

$ws.Range("D45").Value = "'13.42"
$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").Value = "'3.639"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("D48").Value = "'1.210"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").Value = "'122.08"
$ws.Range("E49").Value = "  -2.17%  "

$ws.Range("D50").Value = "'79.03"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'1.127"
$ws.Range("E51").Value = "  -2.58%  "
